$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match rows scraped by the automated tracker sync (event_id, fecha, jugador_A, jugador_B, pronostico, cuota)
$newRows = @(
    @("14494909", "2025-08-27", "Brandon Nakashima", "Jerome Kym", "Gana Jerome Kym", 3.75),
    @("14494908", "2025-08-27", "Jiri Lehecka", "Tomas Martin Etcheverry", "Gana Tomas Martin Etcheverry", 4),
    @("14494917", "2025-08-27", "Joao Fonseca", "Tomas Machac", "Gana Tomas Machac", 1.8),
    @("14494894", "2025-08-27", "Cameron Norrie", "Francisco Comesaña", "Gana Francisco Comesaña", 2.75),
    @("14494907", "2025-08-27", "Jordan Thompson", "Adrian Mannarino", "Gana Jordan Thompson", 2.5),
    @("14495010", "2025-08-27", "Moyuka Uchijima", "Barbora Krejcikova", "Gana Moyuka Uchijima", 5),
    @("14495025", "2025-08-27", "Elena Rybakina", "Tereza Valentova", "Gana Tereza Valentova", 4),
    @("14495018", "2025-08-27", "Elsa Jacquemot", "Leylah Fernandez", "Gana Elsa Jacquemot", 4),
    @("14487481", "2025-08-28", "Luka Mikrut", "Carlo Alberto Caniato", "Gana Carlo Alberto Caniato", 2),
    @("14487488", "2025-08-28", "Nicolai Budkov Kjaer", "Giovanni Fonio", "Gana Giovanni Fonio", 2.5),
    @("14534342", "2025-08-27", "Abdullah Shelbayh", "Daniel Rincon", "Gana Abdullah Shelbayh", 2.25),
    @("14487515", "2025-08-28", "David Jorda Sanchis", "Mika Brunold", "Gana David Jorda Sanchis", 3.4),
    @("14487517", "2025-08-28", "Dimitar Kuzmanov", "Daniel Michalski", "Gana Daniel Michalski", 2.2)
)

$startRow = 474
$r = $startRow
foreach ($row in $newRows) {
    $colA = $ws.Cells.Item($r, 1)
    $colA.NumberFormat = "@"
    $colA.Value = $row[0]

    $colB = $ws.Cells.Item($r, 2)
    $colB.NumberFormat = "@"
    $colB.Value = $row[1]

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = [double]$row[5]

    # Match still pending -> resultado/profit left blank, same as the other upcoming fixtures
    $ws.Cells.Item($r, 7).Value = ""
    $ws.Cells.Item($r, 8).Value = ""

    $r = $r + 1
}
